$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 968.75
$ws.Range("I6").Value = 350
$ws.Range("K6").Value = 1050
$ws.Range("M6").Value = -938
$ws.Range("H38").Value = 1338.6364
$ws.Range("I38").Value = 1883
$ws.Range("J38").Value = 685.4
$ws.Range("K38").Value = 5649
$ws.Range("L38").Value = 2056.2
$ws.Range("M38").Value = -5277
$ws.Range("N38").Value = -2800.2
$ws.Range("H42").Value = 1493.5
$ws.Range("I42").Value = 2289.8
$ws.Range("J42").Value = 166.33333
$ws.Range("K42").Value = 6869.400000000001
$ws.Range("L42").Value = 498.99999
$ws.Range("M42").Value = -6639.400000000001
$ws.Range("N42").Value = -958.99999
$ws.Range("H48").Value = 500
$ws.Range("I48").Value = 500
$ws.Range("K48").Value = 1500
$ws.Range("M48").Value = -1208
$ws.Range("H56").Value = 500
$ws.Range("I56").Value = 500
$ws.Range("K56").Value = 1500
$ws.Range("M56").Value = -966
$ws.Range("H116").Value = 3087.7778
$ws.Range("I116").Value = 3287.5
$ws.Range("J116").Value = 1490
$ws.Range("K116").Value = 3287.5
$ws.Range("L116").Value = 1490
$ws.Range("M116").Value = 154.5
$ws.Range("N116").Value = -8374
$ws.Range("H137").Value = 22223368
$ws.Range("I137").Value = 27778592
$ws.Range("K137").Value = 83335776
$ws.Range("M137").Value = -83333226
$ws.Range("H141").Value = 3601.1936
$ws.Range("I141").Value = 2513.25
$ws.Range("J141").Value = 7331.2856
$ws.Range("K141").Value = 7539.75
$ws.Range("L141").Value = 21993.8568
$ws.Range("M141").Value = -2359.75
$ws.Range("N141").Value = -32353.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5377.67
$ws.Range("I32").Value = 3265.53
$ws.Range("J32").Value = 15689.883
$ws.Range("K32").Value = 3265.53
$ws.Range("L32").Value = 15689.883
$ws.Range("M32").Value = -2978.53
$ws.Range("N32").Value = -16263.883
$ws.Range("H113").Value = 500398
$ws.Range("J113").Value = 500398
$ws.Range("L113").Value = 500398
$ws.Range("N113").Value = -509076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1599.375
$ws.Range("I86").Value = 1565.8334
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 1565.8334
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -442.8334
$ws.Range("N86").Value = -3946
$ws.Range("H89").Value = 1599.375
$ws.Range("I89").Value = 1565.8334
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 7829.166999999999
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -2213.166999999999
$ws.Range("N89").Value = -19732
$ws.Range("H134").Value = 3870.7144
$ws.Range("I134").Value = 3553.077
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 10659.231
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -8124.231
$ws.Range("N134").Value = -29070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6225.85
$ws.Range("I31").Value = 2085.9614
$ws.Range("J31").Value = 13914.214
$ws.Range("K31").Value = 2085.9614
$ws.Range("L31").Value = 13914.214
$ws.Range("M31").Value = -1790.9614
$ws.Range("N31").Value = -14504.214
$ws.Range("H33").Value = 26222.166
$ws.Range("I33").Value = 26222.166
$ws.Range("K33").Value = 26222.166
$ws.Range("M33").Value = -25843.166
$ws.Range("H34").Value = 6225.85
$ws.Range("I34").Value = 2085.9614
$ws.Range("J34").Value = 13914.214
$ws.Range("K34").Value = 2085.9614
$ws.Range("L34").Value = 13914.214
$ws.Range("M34").Value = -1883.9614
$ws.Range("N34").Value = -14318.214
$ws.Range("H36").Value = 4283.3335
$ws.Range("I36").Value = 2250
$ws.Range("J36").Value = 5300
$ws.Range("K36").Value = 2250
$ws.Range("L36").Value = 5300
$ws.Range("M36").Value = -1862
$ws.Range("N36").Value = -6076
$ws.Range("H40").Value = 4283.3335
$ws.Range("I40").Value = 2250
$ws.Range("J40").Value = 5300
$ws.Range("K40").Value = 2250
$ws.Range("L40").Value = 5300
$ws.Range("M40").Value = -2090
$ws.Range("N40").Value = -5620
$ws.Range("H112").Value = 28000
$ws.Range("J112").Value = 28000
$ws.Range("L112").Value = 28000
$ws.Range("N112").Value = -30954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1399.4445
$ws.Range("I23").Value = 3434
$ws.Range("J23").Value = 382.16666
$ws.Range("K23").Value = 10302
$ws.Range("L23").Value = 1146.49998
$ws.Range("M23").Value = -10067
$ws.Range("N23").Value = -1616.49998
$ws.Range("H34").Value = 2306.3125
$ws.Range("J34").Value = 2159.9333
$ws.Range("L34").Value = 6479.7999
$ws.Range("N34").Value = -6647.7999
$ws.Range("H36").Value = 2000
$ws.Range("J36").Value = 2000
$ws.Range("L36").Value = 6000
$ws.Range("N36").Value = -6338
$ws.Range("H60").Value = 2150.9
$ws.Range("J60").Value = 4193.4
$ws.Range("L60").Value = 12580.2
$ws.Range("N60").Value = -13082.2
$ws.Range("H131").Value = 7753524.5
$ws.Range("I131").Value = 516.6667
$ws.Range("J131").Value = 8335000
$ws.Range("K131").Value = 1550.0001
$ws.Range("L131").Value = 25005000
$ws.Range("M131").Value = 3489.9999
$ws.Range("N131").Value = -25015080
$ws.Range("H137").Value = 9185207
$ws.Range("I137").Value = 12501526
$ws.Range("J137").Value = 341688.34
$ws.Range("K137").Value = 37504578
$ws.Range("L137").Value = 1025065.02
$ws.Range("M137").Value = -37499478
$ws.Range("N137").Value = -1035265.02
$ws.Range("H140").Value = 9087.275
$ws.Range("I140").Value = 15721.429
$ws.Range("J140").Value = 2895.4
$ws.Range("K140").Value = 47164.287
$ws.Range("L140").Value = 8686.200000000001
$ws.Range("M140").Value = -41984.287
$ws.Range("N140").Value = -19046.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 49.833332
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 38.285713
$ws.Range("K2").Value = 66
$ws.Range("L2").Value = 38.285713
$ws.Range("M2").Value = 47
$ws.Range("N2").Value = -264.285713
$ws.Range("H3").Value = 943.6667
$ws.Range("I3").Value = 482.16666
$ws.Range("K3").Value = 482.16666
$ws.Range("M3").Value = -366.16666
$ws.Range("H11").Value = 9875750
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3278
$ws.Range("H13").Value = 1207.7778
$ws.Range("J13").Value = 3465
$ws.Range("L13").Value = 3465
$ws.Range("N13").Value = -3743
$ws.Range("H109").Value = 34775
$ws.Range("J109").Value = 34775
$ws.Range("L109").Value = 34775
$ws.Range("N109").Value = -36855
$ws.Range("H113").Value = 1260
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -5840
$ws.Range("H126").Value = 2396.9355
$ws.Range("I126").Value = 1850.2858
$ws.Range("J126").Value = 2847.1177
$ws.Range("K126").Value = 5550.857400000001
$ws.Range("L126").Value = 8541.3531
$ws.Range("M126").Value = -3080.857400000001
$ws.Range("N126").Value = -13481.3531
$ws.Range("H132").Value = 3629.3044
$ws.Range("I132").Value = 2965.9048
$ws.Range("J132").Value = 10595
$ws.Range("K132").Value = 8897.714399999999
$ws.Range("L132").Value = 31785
$ws.Range("M132").Value = -6367.714399999999
$ws.Range("N132").Value = -36845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 683.5714
$ws.Range("I46").Value = 630
$ws.Range("J46").Value = 780
$ws.Range("K46").Value = 630
$ws.Range("L46").Value = 780
$ws.Range("M46").Value = -442
$ws.Range("N46").Value = -1156
$ws.Range("H55").Value = 435
$ws.Range("I55").Value = 288.66666
$ws.Range("J55").Value = 532.55554
$ws.Range("K55").Value = 288.66666
$ws.Range("L55").Value = 532.55554
$ws.Range("M55").Value = -115.66666
$ws.Range("N55").Value = -878.55554
$ws.Range("H106").Value = 20536
$ws.Range("J106").Value = 20536
$ws.Range("L106").Value = 20536
$ws.Range("N106").Value = -23060
$ws.Range("H132").Value = 4021.054
$ws.Range("I132").Value = 2771.1365
$ws.Range("J132").Value = 5854.2666
$ws.Range("K132").Value = 8313.4095
$ws.Range("L132").Value = 17562.7998
$ws.Range("M132").Value = -5783.4095
$ws.Range("N132").Value = -22622.7998
$ws.Range("H136").Value = 4191.2964
$ws.Range("I136").Value = 2242.2856
$ws.Range("J136").Value = 11012.833
$ws.Range("K136").Value = 6726.8568
$ws.Range("L136").Value = 33038.499
$ws.Range("M136").Value = -4176.8568
$ws.Range("N136").Value = -38138.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 29926.666
$ws.Range("J97").Value = 29926.666
$ws.Range("L97").Value = 29926.666
$ws.Range("N97").Value = -31908.666
